$wb = $excel.ActiveWorkbook

# --- "Overview" sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G2").Value = "2016-08-15 22:15:39"
$ws.Range("G4").Value = "2016-08-15 22:15:39"

# --- "zh-cn" sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E2").Value = "mt"
$ws.Range("E4").Value = "mt"
$ws.Range("H2").Value = "2016-08-15 22:15:34"
$ws.Range("H4").Value = "2016-08-15 22:15:34"
$ws.Range("K2").Value = "2016-08-15 22:15:51"
$ws.Range("K4").Value = "2016-08-15 22:15:51"

# --- "de-de" sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("H2").Value = "2016-08-15 22:15:39"
$ws.Range("H4").Value = "2016-08-15 22:15:39"
$ws.Range("K2").Value = "2016-08-15 22:15:58"
$ws.Range("K4").Value = "2016-08-15 22:15:58"
